$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.046.27'
$ws.Range("E2").Value = '  -3.67%  '
$ws.Range("D3").Value = '2.197.99'
$ws.Range("E3").Value = '  -3.75%  '
$ws.Range("E4").Value = '  +0.49%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '105.55'
$ws.Range("E5").Value = '  -15.12%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '295.36'
$ws.Range("E6").Value = '  +10.63%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.617'
$ws.Range("E7").Value = '  -3.71%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.587'
$ws.Range("E9").Value = '  -6.22%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '42.87'
$ws.Range("E10").Value = '  -11.26%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0901'
$ws.Range("E11").Value = '  -5.48%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '53.92'
$ws.Range("E12").Value = '  -0.78%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '8.61'
$ws.Range("E13").Value = '  -7.96%  '
$ws.Range("E14").Value = '  -3.68%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.923'
$ws.Range("E15").Value = '  +1.48%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.76'
$ws.Range("E16").Value = '  -4.97%  '
$ws.Range("D17").Value = '2.534.08'
$ws.Range("E17").Value = '  -3.34%  '
$ws.Range("D18").Value = '2.233.38'
$ws.Range("E18").Value = '  -2.03%  '
$ws.Range("D19").Value = '41.867.84'
$ws.Range("E19").Value = '  -4.22%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.26'
$ws.Range("E20").Value = '  +3.80%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0000103'
$ws.Range("E21").Value = '  -6.71%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '71.91'
$ws.Range("E22").Value = '  -0.71%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.45'
$ws.Range("E23").Value = '  +19.05%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.24'
$ws.Range("E24").Value = '  -8.11%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '225.92'
$ws.Range("E25").Value = '  -4.21%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.79'
$ws.Range("E26").Value = '  -7.68%  '
$ws.Range("E27").Value = '  -1.62%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '11.36'
$ws.Range("E28").Value = '  -5.18%  '
$ws.Range("E29").Value = '  +0.04%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.23'
$ws.Range("E30").Value = '  -1.33%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '37.49'
$ws.Range("E31").Value = '  -11.77%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.18'
$ws.Range("E32").Value = '  -5.33%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '172.30'
$ws.Range("E33").Value = '  -0.46%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '20.65'
$ws.Range("E34").Value = '  -5.06%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0865'
$ws.Range("E35").Value = '  -6.56%  '
$ws.Range("B36").Value = 'Filecoin'
$ws.Range("C36").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.44'
$ws.Range("E36").Value = '  -6.00%  '
$ws.Range("B37").Value = 'RenderToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.88'
$ws.Range("E37").Value = '  +5.16%  '
$ws.Range("E38").Value = '  -2.60%  '
$ws.Range("E39").Value = '  -4.52%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0357'
$ws.Range("E40").Value = '  -5.54%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.100'
$ws.Range("E41").Value = '  -5.98%  '
$ws.Range("E42").Value = '  -5.03%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '69.64'
$ws.Range("E43").Value = '  -6.15%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.225'
$ws.Range("E44").Value = '  -5.75%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.43'
$ws.Range("E46").Value = '  -11.41%  '
$ws.Range("E47").Value = '  -7.32%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.34'
$ws.Range("E48").Value = '  -5.63%  '
$ws.Range("E49").Value = '  +2.55%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '101.31'
$ws.Range("E50").Value = '  -1.00%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.36'
$ws.Range("E51").Value = '  -2.98%  '
